# Auto-update draw results: append the 2025-10-13 Pick 3 draw as a new
# row at the bottom of the Results sheet (row 27), matching the nightly
# scraper's existing layout: Date | Game | Phase | Result | InsertedAt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 27

# Columns A (Date) and C (Phase) look like a date / a plain number
# ("2025-10-13", "251013") and would otherwise be auto-coerced to a
# numeric/date value on assignment. Force them to literal Text first
# (same intent as typing them into a pre-"Format Cells > Text" column),
# exactly like every other row already stored in this sheet.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("C" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-10-13"
$ws.Range("B" + $newRow).Value = "Pick 3"
$ws.Range("C" + $newRow).Value = "251013"
$ws.Range("D" + $newRow).Value = "6-6-2"
$ws.Range("E" + $newRow).Value = "2025-10-13T21:37:36.299+04:00"

# The sheet suppresses the "number stored as text" checker across the
# whole data range; keep that marker extended onto the newly-added row.
$ws.Range("A1:E" + $newRow).Errors.Item(1).Ignore = $true
